$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-68 (A:E), replacing the previous A1:E69 range
$data = @(
    @(64, 847, 17, 77, 86),
    @(61, 847, 131, 79, 107),
    @(54, 852, 259, 61, 104),
    @(42, 852, 404, 67, 11),
    @(41, 853, 429, 65, 11),
    @(40, 855, 459, 66, 12),
    @(67, 682, 0, 88, 119),
    @(58, 682, 140, 80, 100),
    @(52, 676, 263, 90, 96),
    @(46, 685, 378, 74, 116),
    @(36, 685, 517, 59, 92),
    @(30, 685, 643, 74, 78),
    @(28, 685, 756, 76, 100),
    @(23, 679, 882, 82, 111),
    @(17, 684, 1022, 75, 83),
    @(10, 685, 1140, 70, 90),
    @(5, 687, 1265, 61, 97),
    @(1, 686, 1393, 73, 84),
    @(65, 523, 10, 72, 111),
    @(57, 527, 142, 70, 93),
    @(55, 516, 254, 81, 110),
    @(43, 523, 393, 76, 87),
    @(35, 520, 520, 84, 79),
    @(34, 523, 620, 58, 129),
    @(26, 519, 766, 66, 95),
    @(18, 521, 1013, 68, 11),
    @(16, 523, 1043, 63, 10),
    @(15, 523, 1072, 67, 13),
    @(14, 516, 1116, 76, 119),
    @(8, 518, 1250, 59, 104),
    @(3, 520, 1376, 46, 93),
    @(66, 368, 10, 78, 116),
    @(59, 355, 137, 89, 100),
    @(51, 361, 265, 76, 103),
    @(45, 370, 389, 65, 81),
    @(39, 364, 502, 67, 94),
    @(33, 366, 625, 61, 108),
    @(29, 360, 752, 80, 113),
    @(22, 359, 890, 79, 98),
    @(20, 362, 1006, 59, 99),
    @(11, 358, 1128, 66, 94),
    @(9, 357, 1240, 80, 113),
    @(4, 355, 1362, 66, 107),
    @(62, 190, 21, 80, 86),
    @(56, 193, 145, 75, 90),
    @(50, 194, 270, 69, 10),
    @(49, 197, 300, 65, 10),
    @(48, 197, 333, 68, 10),
    @(47, 190, 377, 80, 108),
    @(37, 189, 512, 75, 87),
    @(31, 191, 637, 73, 91),
    @(25, 189, 879, 73, 100),
    @(19, 198, 1013, 65, 83),
    @(12, 205, 1128, 39, 95),
    @(7, 190, 1255, 76, 90),
    @(2, 193, 1380, 74, 78),
    @(63, 31, 19, 82, 83),
    @(60, 27, 136, 88, 99),
    @(53, 27, 261, 87, 110),
    @(44, 34, 392, 87, 86),
    @(38, 41, 509, 42, 92),
    @(32, 35, 633, 73, 96),
    @(27, 33, 762, 75, 84),
    @(24, 29, 882, 79, 84),
    @(21, 22, 1001, 85, 92),
    @(13, 25, 1120, 81, 117),
    @(6, 32, 1259, 88, 82)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
    }
}

# Remove the now-obsolete last row (69) entirely, shrinking the used range to A1:E68
$ws.Rows.Item(69).Delete()

$wb.Save()
